# Auto-generated edit script: updates crypto price/volume table to match
# the target OOXML diff (commit: "Updated cryptos list ... GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.714.33"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.633.78"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'217.97"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "1.861.38"
$ws.Range("D13").Value = "1.627.48"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("E14").Value = "  -2.79%  "
$ws.Range("D15").Value = "'0.521"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "'64.04"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").Value = "26.690.73"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("E18").Value = "  -3.12%  "
$ws.Range("D19").Value = "'211.28"
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("D23").Value = "'2.33"
$ws.Range("E23").Value = "  -4.14%  "
$ws.Range("E24").Value = "  -3.22%  "
$ws.Range("D25").Value = "'146.83"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("E28").Value = "  -3.14%  "
$ws.Range("D29").Value = "'15.54"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("E30").Value = "  -4.53%  "
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").Value = "'2.94"
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("D34").Value = "1.261.47"
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("E37").Value = "  -3.48%  "
$ws.Range("D38").Value = "'0.524"
$ws.Range("E38").Value = "  -3.60%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "'0.802"
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("D41").Value = "'0.798"
$ws.Range("E41").Value = "  -2.66%  "
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("E43").Value = "  -4.48%  "
$ws.Range("D44").Value = "1.771.84"
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("D45").Value = "'91.44"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").Value = "'59.67"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  -4.04%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.407"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0957"
$ws.Range("E51").Value = "  -2.82%  "
